$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 4
$ws.Range("C9").Value = 30
$ws.Range("C10").Value = 316
$ws.Range("C11").Value = 3206
$ws.Range("C12").Value = 32240
$ws.Range("C13").Value = 9
$ws.Range("C14").Value = 83
$ws.Range("C15").Value = 797
$ws.Range("C16").Value = 8188
$ws.Range("C17").Value = 80997

$excel.ActiveWindow.Zoom = 164
$ws.Range("C18").Select()
